$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 22.88430570846026
$ws.Cells.Item(2, 3).Value = 11.38541711786571
$ws.Cells.Item(2, 5).Value = 8.594639669030494
$ws.Cells.Item(2, 6).Value = 16.86991607391233
$ws.Cells.Item(2, 7).Value = 3.703260823887519
$ws.Cells.Item(2, 9).Value = 30.67623233654734
$ws.Cells.Item(2, 12).Value = 10.73578321317222
$ws.Cells.Item(2, 14).Value = 19.48541697895221

$ws.Cells.Item(3, 2).Value = 22.42047879177176
$ws.Cells.Item(3, 3).Value = 10.83131493425334
$ws.Cells.Item(3, 5).Value = 8.603163084067516
$ws.Cells.Item(3, 6).Value = 15.89584955866808
$ws.Cells.Item(3, 7).Value = 3.707719094281989
$ws.Cells.Item(3, 9).Value = 30.67206585270842
$ws.Cells.Item(3, 12).Value = 10.71740860601133
$ws.Cells.Item(3, 14).Value = 19.56045615136209

$ws.Cells.Item(4, 2).Value = 22.13809013655271
$ws.Cells.Item(4, 3).Value = 10.4800525084757
$ws.Cells.Item(4, 5).Value = 8.608812154364172
$ws.Cells.Item(4, 6).Value = 15.26997757108491
$ws.Cells.Item(4, 7).Value = 3.710594896917085
$ws.Cells.Item(4, 9).Value = 30.67780555932537
$ws.Cells.Item(4, 12).Value = 10.70848734333051
$ws.Cells.Item(4, 14).Value = 19.60859380239062

$ws.Cells.Item(5, 2).Value = 22.02379082481288
$ws.Cells.Item(5, 3).Value = 10.33435132713477
$ws.Cells.Item(5, 5).Value = 8.611218856521949
$ws.Cells.Item(5, 6).Value = 15.00819731993403
$ws.Cells.Item(5, 7).Value = 3.711801759845825
$ws.Cells.Item(5, 9).Value = 30.68222199808178
$ws.Cells.Item(5, 12).Value = 10.70544700305839
$ws.Cells.Item(5, 14).Value = 19.62873025964426

$ws.Cells.Item(6, 2).Value = 22.00486345471049
$ws.Cells.Item(6, 3).Value = 10.31001019978839
$ws.Cells.Item(6, 5).Value = 8.611624812993275
$ws.Cells.Item(6, 6).Value = 14.96433081551593
$ws.Cells.Item(6, 7).Value = 3.712004273706056
$ws.Cells.Item(6, 9).Value = 30.68308052483995
$ws.Cells.Item(6, 12).Value = 10.70497813666233
$ws.Cells.Item(6, 14).Value = 19.63210534189693

$ws.Cells.Item(7, 2).Value = 22.13654528422283
$ws.Cells.Item(7, 3).Value = 10.47809758105379
$ws.Cells.Item(7, 5).Value = 8.608844188034425
$ws.Cells.Item(7, 6).Value = 15.26647399323728
$ws.Cells.Item(7, 7).Value = 3.710611031385795
$ws.Cells.Item(7, 9).Value = 30.67785672238547
$ws.Cells.Item(7, 12).Value = 10.70844392887024
$ws.Cells.Item(7, 14).Value = 19.6088632627947

$ws.Cells.Item(8, 2).Value = 22.7239821956056
$ws.Cells.Item(8, 3).Value = 11.19677513814738
$ws.Cells.Item(8, 5).Value = 8.597492354979149
$ws.Cells.Item(8, 6).Value = 16.53996406344769
$ws.Cells.Item(8, 7).Value = 3.704769397691829
$ws.Cells.Item(8, 9).Value = 30.67306840155383
$ws.Cells.Item(8, 12).Value = 10.72895869886002
$ws.Cells.Item(8, 14).Value = 19.51086293810517

$ws.Cells.Item(9, 2).Value = 23.88756609412613
$ws.Cells.Item(9, 3).Value = 12.51047203485275
$ws.Cells.Item(9, 5).Value = 8.578523572486265
$ws.Cells.Item(9, 6).Value = 19.0027458068253
$ws.Cells.Item(9, 7).Value = 3.694405472344733
$ws.Cells.Item(9, 9).Value = 30.72986919414148
$ws.Cells.Item(9, 12).Value = 10.78783478117294
$ws.Cells.Item(9, 14).Value = 19.33500233240787

$ws.Cells.Item(10, 2).Value = 24.7398573562466
$ws.Cells.Item(10, 3).Value = 13.4085093012062
$ws.Cells.Item(10, 5).Value = 8.566586474858132
$ws.Cells.Item(10, 6).Value = 20.67494806633232
$ws.Cells.Item(10, 7).Value = 3.687447121904872
$ws.Cells.Item(10, 9).Value = 30.81234615036747
$ws.Cells.Item(10, 12).Value = 10.8423160469218
$ws.Cells.Item(10, 14).Value = 19.21566819853621

$ws.Cells.Item(11, 2).Value = 25.12509115630331
$ws.Cells.Item(11, 3).Value = 13.80107222504272
$ws.Cells.Item(11, 5).Value = 8.561588590695978
$ws.Cells.Item(11, 6).Value = 21.3917225636224
$ws.Cells.Item(11, 7).Value = 3.684422025453419
$ws.Cells.Item(11, 9).Value = 30.85876838362413
$ws.Cells.Item(11, 12).Value = 10.86949819437495
$ws.Cells.Item(11, 14).Value = 19.16350715867067

$ws.Cells.Item(12, 2).Value = 25.27046196458691
$ws.Cells.Item(12, 3).Value = 13.94733959163762
$ws.Cells.Item(12, 5).Value = 8.559758087330994
$ws.Cells.Item(12, 6).Value = 21.65686569030329
$ws.Cells.Item(12, 7).Value = 3.683296518145761
$ws.Cells.Item(12, 9).Value = 30.87762925568592
$ws.Cells.Item(12, 12).Value = 10.88013204117402
$ws.Cells.Item(12, 14).Value = 19.14405954428281

$ws.Cells.Item(13, 2).Value = 25.23917873011668
$ws.Cells.Item(13, 3).Value = 13.91594579755452
$ws.Cells.Item(13, 5).Value = 8.560149558964698
$ws.Cells.Item(13, 6).Value = 21.60004134736742
$ws.Cells.Item(13, 7).Value = 3.683538027666532
$ws.Cells.Item(13, 9).Value = 30.8735102140651
$ws.Cells.Item(13, 12).Value = 10.87782677926352
$ws.Cells.Item(13, 14).Value = 19.14823440079219

$ws.Cells.Item(14, 2).Value = 25.1370619421678
$ws.Cells.Item(14, 3).Value = 13.81315405649401
$ws.Cells.Item(14, 5).Value = 8.561436750237949
$ws.Cells.Item(14, 6).Value = 21.4136618050453
$ws.Cells.Item(14, 7).Value = 3.684329028654786
$ws.Cells.Item(14, 9).Value = 30.86029438693107
$ws.Cells.Item(14, 12).Value = 10.87036624844276
$ws.Cells.Item(14, 14).Value = 19.1619010912614

$ws.Cells.Item(15, 2).Value = 25.07444180876263
$ws.Cells.Item(15, 3).Value = 13.74987772132286
$ws.Cells.Item(15, 5).Value = 8.562233275160466
$ws.Cells.Item(15, 6).Value = 21.29868154950795
$ws.Cells.Item(15, 7).Value = 3.684816144000513
$ws.Cells.Item(15, 9).Value = 30.85236625556957
$ws.Cells.Item(15, 12).Value = 10.86584067249367
$ws.Cells.Item(15, 14).Value = 19.17031198111875

$ws.Cells.Item(16, 2).Value = 24.7146183268018
$ws.Cells.Item(16, 3).Value = 13.38252494528328
$ws.Cells.Item(16, 5).Value = 8.566921791055847
$ws.Cells.Item(16, 6).Value = 20.62722412089977
$ws.Cells.Item(16, 7).Value = 3.687647629899992
$ws.Cells.Item(16, 9).Value = 30.80949166697436
$ws.Cells.Item(16, 12).Value = 10.84058753154429
$ws.Cells.Item(16, 14).Value = 19.21911972242866

$ws.Cells.Item(17, 2).Value = 24.49313243817534
$ws.Cells.Item(17, 3).Value = 13.15300710078527
$ws.Cells.Item(17, 5).Value = 8.569908718594462
$ws.Cells.Item(17, 6).Value = 20.20408069617459
$ws.Cells.Item(17, 7).Value = 3.68942048532094
$ws.Cells.Item(17, 9).Value = 30.78547136568109
$ws.Cells.Item(17, 12).Value = 10.8257068926957
$ws.Cells.Item(17, 14).Value = 19.24960522764432

$ws.Cells.Item(18, 2).Value = 24.36551500376209
$ws.Cells.Item(18, 3).Value = 13.0194968047054
$ws.Cells.Item(18, 5).Value = 8.571667420989302
$ws.Cells.Item(18, 6).Value = 19.95656407809808
$ws.Cells.Item(18, 7).Value = 3.690453398350939
$ws.Cells.Item(18, 9).Value = 30.77249352891575
$ws.Cells.Item(18, 12).Value = 10.81737397487369
$ws.Cells.Item(18, 14).Value = 19.26733971762741

$ws.Cells.Item(19, 2).Value = 24.32227214333912
$ws.Cells.Item(19, 3).Value = 12.97403838318318
$ws.Cells.Item(19, 5).Value = 8.572269880513629
$ws.Cells.Item(19, 6).Value = 19.87204792380562
$ws.Cells.Item(19, 7).Value = 3.690805398733119
$ws.Cells.Item(19, 9).Value = 30.76824331219118
$ws.Cells.Item(19, 12).Value = 10.81459153744314
$ws.Cells.Item(19, 14).Value = 19.27337869456154

$ws.Cells.Item(20, 2).Value = 24.51673436754653
$ws.Cells.Item(20, 3).Value = 13.1775953984219
$ws.Cells.Item(20, 5).Value = 8.569586543337266
$ws.Cells.Item(20, 6).Value = 20.2495528364879
$ws.Cells.Item(20, 7).Value = 3.689230395324453
$ws.Cells.Item(20, 9).Value = 30.78794162137691
$ws.Cells.Item(20, 12).Value = 10.82726759850187
$ws.Cells.Item(20, 14).Value = 19.24633929735316

$ws.Cells.Item(21, 2).Value = 25.16707110215989
$ws.Cells.Item(21, 3).Value = 13.84341194352456
$ws.Cells.Item(21, 5).Value = 8.561056986595831
$ws.Cells.Item(21, 6).Value = 21.46857628470567
$ws.Cells.Item(21, 7).Value = 3.684096149989625
$ws.Cells.Item(21, 9).Value = 30.86414140409537
$ws.Cells.Item(21, 12).Value = 10.87254838125571
$ws.Cells.Item(21, 14).Value = 19.15787859204071

$ws.Cells.Item(22, 2).Value = 25.58907304487303
$ws.Cells.Item(22, 3).Value = 14.26461214278993
$ws.Cells.Item(22, 5).Value = 8.55584428263937
$ws.Cells.Item(22, 6).Value = 22.22866616901555
$ws.Cells.Item(22, 7).Value = 3.680857314076856
$ws.Cells.Item(22, 9).Value = 30.92141317005181
$ws.Cells.Item(22, 12).Value = 10.90412469812984
$ws.Cells.Item(22, 14).Value = 19.10183969098501

$ws.Cells.Item(23, 2).Value = 25.36416814282596
$ws.Cells.Item(23, 3).Value = 14.04111262929909
$ws.Cells.Item(23, 5).Value = 8.55859331798564
$ws.Cells.Item(23, 6).Value = 21.82633154475864
$ws.Cells.Item(23, 7).Value = 3.682575311729361
$ws.Cells.Item(23, 9).Value = 30.89016246856743
$ws.Cells.Item(23, 12).Value = 10.88709195004303
$ws.Cells.Item(23, 14).Value = 19.13158655228174

$ws.Cells.Item(24, 2).Value = 24.50606480274365
$ws.Cells.Item(24, 3).Value = 13.16648387289493
$ws.Cells.Item(24, 5).Value = 8.569732069569753
$ws.Cells.Item(24, 6).Value = 20.22900810905294
$ws.Cells.Item(24, 7).Value = 3.689316292406653
$ws.Cells.Item(24, 9).Value = 30.78682222818184
$ws.Cells.Item(24, 12).Value = 10.82656131110657
$ws.Cells.Item(24, 14).Value = 19.24781517660005

$ws.Cells.Item(25, 2).Value = 23.57256191358578
$ws.Cells.Item(25, 3).Value = 12.1663278047407
$ws.Cells.Item(25, 5).Value = 8.583303503259021
$ws.Cells.Item(25, 6).Value = 18.34778573295697
$ws.Cells.Item(25, 7).Value = 3.697093306777578
$ws.Cells.Item(25, 9).Value = 30.70737042208741
$ws.Cells.Item(25, 12).Value = 10.7699237731589
$ws.Cells.Item(25, 14).Value = 19.38083842887077
